# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the per-row observations for this
# market/product subset: rows 2,4,5,10 cycle among themselves, rows 6/7
# swap, and rows 9/11 swap. Row 3 and row 8 are untouched by this update.
# Apply it by writing each row's new Fecha / Volumen / Precio* / Origen
# values directly (rather than relying on any sort), so every cell ends
# up exactly as the refreshed source data dictates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44214
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 1800
$ws.Range("O2").Value = 1800
$ws.Range("P2").Value = 1800
$ws.Range("S2").Value = 1800

# Row 4
$ws.Range("D4").Value = 44592
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = 7500
$ws.Range("O4").Value = 7500
$ws.Range("P4").Value = 7500
$ws.Range("S4").Value = 7500

# Row 5
$ws.Range("D5").Value = 44323
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 3200
$ws.Range("O5").Value = 3200
$ws.Range("P5").Value = 3200
$ws.Range("S5").Value = 3200

# Row 6
$ws.Range("D6").Value = 44176
$ws.Range("M6").Value = 20
$ws.Range("R6").Value = "Región de O'Higgins"

# Row 7
$ws.Range("D7").Value = 44574
$ws.Range("M7").Value = 200
$ws.Range("R7").Value = "Región de La Araucanía"

# Row 9
$ws.Range("D9").Value = 44567
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 2400
$ws.Range("O9").Value = 2400
$ws.Range("P9").Value = 2400
$ws.Range("R9").Value = "Región de La Araucanía"
$ws.Range("S9").Value = 2400

# Row 10
$ws.Range("D10").Value = 44616
$ws.Range("M10").Value = 200

# Row 11
$ws.Range("D11").Value = 44551
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 4500
$ws.Range("O11").Value = 4500
$ws.Range("P11").Value = 4500
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 4500
